# update new orleans xlsx files
#
# The two worksheets swap identities: the physical sheet that used to be
# "hotel_info" becomes "review_info" (headers only), and the physical sheet
# that used to be "review_info" becomes "hotel_info" (original hotel data
# plus a new "State" column inserted right after "Hotel_Name").

$wb = $excel.ActiveWorkbook

$sheetA = $wb.Worksheets.Item("hotel_info")   # currently hotel_info
$sheetB = $wb.Worksheets.Item("review_info")  # currently review_info

# --- swap the two sheet names without colliding -----------------------
$sheetA.Name = "__tmp_swap__"
$sheetB.Name = "hotel_info"
$sheetA.Name = "review_info"

# $sheetA is now named "review_info" -> give it the review_info headers
# $sheetB is now named "hotel_info"  -> give it the hotel_info data (+ State)

$reviewSheet = $sheetA
$hotelSheet = $sheetB

# --- wipe existing contents --------------------------------------------
$reviewSheet.Cells.Clear()
$hotelSheet.Cells.Clear()

# --- review_info: header row only --------------------------------------
$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $reviewSheet.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- hotel_info: header row with new "State" column --------------------
$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $hotelSheet.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

# --- hotel_info: data row, State="Louisiana" inserted after Hotel_Name -
$hotelSheet.Cells.Item(2, 1).Value = 57568
$hotelSheet.Cells.Item(2, 2).Value = "Suburban Extended Stay Hotel La Place"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"
$hotelSheet.Cells.Item(2, 4).Value = "La Place"
$hotelSheet.Cells.Item(2, 5).Value = 70068
$hotelSheet.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g40267-d1231096-Reviews-Suburban_Extended_Stay-LaPlace_Louisiana.html"
$hotelSheet.Cells.Item(2, 7).Value = "Suburban Extended Stay"

# H2:J2 hold digit-only strings ("47", "7", "47") that must stay TEXT
# (matching the source data), not get auto-coerced to numbers.
$hotelSheet.Range("H2:J2").NumberFormat = "@"
$hotelSheet.Cells.Item(2, 8).Value = "47"
$hotelSheet.Cells.Item(2, 9).Value = "7"
$hotelSheet.Cells.Item(2, 10).Value = "47"
